$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This is the weekly data refresh: a new observation (dated 2022-08-25,
# serial 44798) is inserted right after the existing row 62, pushing the
# previously-recorded rows 63-75 down to 64-76 (last of which duplicates
# the old row 75 content).
$ws.Rows.Item(63).Insert()

$newRow = 63
$ws.Cells.Item($newRow, 1).Value2  = 1
$ws.Cells.Item($newRow, 2).Value2  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item($newRow, 3).Value2  = "Arica y Parinacota"
$ws.Cells.Item($newRow, 4).Value2  = 44798
$ws.Cells.Item($newRow, 5).Value2  = 15
$ws.Cells.Item($newRow, 6).Value2  = 100112012
$ws.Cells.Item($newRow, 7).Value2  = "Espinaca"
$ws.Cells.Item($newRow, 8).Value2  = "Sin especificar"
$ws.Cells.Item($newRow, 9).Value2  = "Primera"
$ws.Cells.Item($newRow, 10).Value2 = 300
$ws.Cells.Item($newRow, 11).Value2 = 1800
$ws.Cells.Item($newRow, 12).Value2 = 2000
$ws.Cells.Item($newRow, 13).Value2 = 1900
$ws.Cells.Item($newRow, 14).Value2 = "`$/atado 2,5 a 3 kilos"
$ws.Cells.Item($newRow, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item($newRow, 16).Value2 = 633
$ws.Cells.Item($newRow, 17).Value2 = 3
$ws.Cells.Item($newRow, 18).Value2 = "Hortaliza"
